# Updates the "Price" (column D) and "Volume(1h)" (column E) figures for the
# cryptos list, plus swaps the PaxDollar / TrustWalletToken rows (44 & 45)
# to reflect their new ranking order with refreshed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceCell {
    param(
        [int]$Row,
        $Price,
        [string]$Volume
    )

    if ($null -ne $Price) {
        # Prefix with an apostrophe so Excel stores the figure as literal text
        # (preserves thousand-separator dots / trailing zeros instead of
        # coercing to a number), then strip the resulting "Text" number
        # format back to the sheet's default style.
        $ws.Cells.Item($Row, 4).Value = "'" + $Price
        $ws.Cells.Item($Row, 4).Style = "Normal"
    }

    if ($null -ne $Volume) {
        $ws.Cells.Item($Row, 5).Value = "  " + $Volume + "  "
    }
}

$updates = @(
    @{ Row = 2;  D = "30.278.43";    E = "-0.06%" },
    @{ Row = 3;  D = "1.929.88";     E = "+0.04%" },
    @{ Row = 4;  D = "1.002";        E = "+0.26%" },
    @{ Row = 5;  D = "249.50";       E = "+0.35%" },
    @{ Row = 6;  D = "0.7166";       E = "-1.59%" },
    @{ Row = 7;  D = "1.002";        E = "+0.25%" },
    @{ Row = 8;  D = "0.3208";       E = "-1.61%" },
    @{ Row = 9;  D = "27.44";        E = "+1.09%" },
    @{ Row = 10; D = "0.07092";      E = "+4.35%" },
    @{ Row = 11; D = "0.7919";       E = "-1.06%" },
    @{ Row = 12; D = "0.08036";      E = "-0.41%" },
    @{ Row = 13; D = "1.927.86";     E = "-0.14%" },
    @{ Row = 14; D = "5.359";        E = "-0.88%" },
    @{ Row = 15; D = "94.68";        E = "+0.20%" },
    @{ Row = 16; D = "14.59";        E = "+0.72%" },
    @{ Row = 17; D = "30.275.77";    E = "-0.06%" },
    @{ Row = 18; D = "256.51";       E = "+0.90%" },
    @{ Row = 19; D = "0.000008052";  E = "+0.73%" },
    @{ Row = 20; D = "5.743";        E = "-1.01%" },
    @{ Row = 21; D = "2.179.61";     E = "-0.26%" },
    @{ Row = 22; D = "1.002";        E = "+0.13%" },
    @{ Row = 23; D = $null;          E = "+0.38%" },
    @{ Row = 24; D = "6.817";        E = "-0.49%" },
    @{ Row = 25; D = "9.533";        E = "-1.34%" },
    @{ Row = 26; D = "165.18";       E = "+3.97%" },
    @{ Row = 27; D = "19.20";        E = "+0.65%" },
    @{ Row = 28; D = "2.290";        E = "-3.02%" },
    @{ Row = 29; D = "0.1278";       E = "-4.79%" },
    @{ Row = 30; D = "1.356";        E = "+0.58%" },
    @{ Row = 31; D = "1.533";        E = "-1.43%" },
    @{ Row = 32; D = "4.401";        E = "+0.39%" },
    @{ Row = 33; D = "4.141";        E = "-0.97%" },
    @{ Row = 34; D = "0.05173";      E = "+2.30%" },
    @{ Row = 35; D = "1.253";        E = "+3.37%" },
    @{ Row = 36; D = "0.7435";       E = "+0.71%" },
    @{ Row = 37; D = "2.767";        E = "+0.53%" },
    @{ Row = 38; D = "0.01952";      E = "-0.70%" },
    @{ Row = 39; D = "2.808";        E = "-0.47%" },
    @{ Row = 40; D = "77.51";        E = "-1.69%" },
    @{ Row = 41; D = "6.356";        E = "-3.43%" },
    @{ Row = 42; D = "0.4493";       E = "+1.07%" },
    @{ Row = 43; D = "1.989";        E = "+0.09%" },
    @{ Row = 46; D = "101.03";       E = "-0.69%" },
    @{ Row = 47; D = "9.748";        E = "+0.57%" },
    @{ Row = 48; D = "7.437";        E = "+2.51%" },
    @{ Row = 49; D = "36.43";        E = "+0.57%" },
    @{ Row = 50; D = "0.06107";      E = "+2.95%" },
    @{ Row = 51; D = "0.4171";       E = "+2.73%" }
)

foreach ($u in $updates) {
    # NOTE: named-parameter binding is unreliable for user-defined functions
    # in this host, so args are passed positionally: Row, Price, Volume.
    Set-PriceCell $u.Row $u.D $u.E
}

# Rows 44/45 swap places: TrustWalletToken now ranks above PaxDollar.
$ws.Cells.Item(44, 2).Value = "TrustWalletToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-PriceCell 44 "0.8418" "+0.97%"

$ws.Cells.Item(45, 2).Value = "PaxDollar"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-PriceCell 45 "1.001" "+0.12%"
